# Automatische test-sync: 2025-08-18 19:45:50
#
# Appends a new log row (row 3) to the "Logs" sheet, extends the
# conditional-formatting ranges that were scoped to row 2 so they also
# cover row 3, and bumps the matching "Dashboard" tally from 1 to 2.

$wb = $excel.ActiveWorkbook

# ---- Logs sheet: append row 3 -------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A3").Value = "Interne taak"
$logs.Range("B3").Value = "kwaliteit@testbedrijf123.nl"
$logs.Range("C3").Value = "Leg dit even neer bij Koen."
$logs.Range("D3").Value = "Intern verzoek / Actie voor medewerker"
$logs.Range("E3").Value = "Bedankt, we hebben dit doorgestuurd naar support@testbedrijf123.nl."
$logs.Range("F3").Value = "2025-08-18 19:45:06"
$logs.Range("G3").Value = "Nee"
$logs.Range("H3").Value = "Ja"
$logs.Range("I3").Value = "Nee"
$logs.Range("J3").Value = "Nee"

# ---- Extend conditional formatting from row 2 to rows 2:3 ----------------
$columns = @("D", "G", "H", "I", "J")
foreach ($col in $columns) {
    $fcs = $logs.Range($col + "2").FormatConditions
    $newRange = $logs.Range($col + "2:" + $col + "3")
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($newRange)
    }
}

# ---- Dashboard sheet: bump the tally for the matching category ----------
$dashboard = $wb.Worksheets.Item("Dashboard")
$dashboard.Range("B2").Value = 2
